# ForwardSampling.pptx edit:
#  - Insert a new "Title and Content" slide at position 11 titled
#    "1) Discussion" with a discussion outline about DataObjects /
#    OutStreams / Steps and Sequence (commit: "7. time dependent sampling").
#  - On the later "A Brief Pause..." slide, change the bullet
#    "Curious behaviors?" to "Curious results?".

$p = $ppt.ActivePresentation

# --- 1) Insert the new "1) Discussion" slide at position 11 -----------------
$layout = $p.SlideMaster.CustomLayouts.Item(2)   # "Title and Content"
$newSlide = $p.Slides.AddSlide(11, $layout)

# Title
$newSlide.Shapes.Item(1).TextFrame.TextRange.Text = "1) Discussion"

# Body content placeholder
$body = $newSlide.Shapes.Item(2).TextFrame.TextRange
$body.Text = "DataObjects: PointSets`rStores one row per sample`rMetadata: ProbabilityWeight, Prefix, PointProbability`r`rOutStreams`rNothing is written to file without an OutStream Print or Plot!`r`rSteps and Sequence`rSteps define possible actions to take`rSequence defines order of actions"

# Paragraph indent levels (COM IndentLevel is 1-based; XML lvl=0 -> IndentLevel 1)
$levels = @(1,2,2,2,1,2,2,1,2,2)
for ($i = 1; $i -le $levels.Length; $i++) {
    $para = $body.Paragraphs($i, 1)
    $para.IndentLevel = $levels[$i - 1]
}

# --- 2) Wording tweak on the "A Brief Pause..." slide ------------------------
for ($i = 1; $i -le $p.Slides.Count; $i++) {
    $slide = $p.Slides.Item($i)
    $titleShape = $slide.Shapes.Item(1)
    if ($titleShape.TextFrame.TextRange.Text -eq "A Brief Pause…") {
        $body2 = $slide.Shapes.Item(2).TextFrame.TextRange
        $n = $body2.Paragraphs().Count
        for ($j = 1; $j -le $n; $j++) {
            $para2 = $body2.Paragraphs($j, 1)
            if ($para2.Text -eq "Curious behaviors?") {
                $para2.Text = "Curious results?"
            }
        }
    }
}
